$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 568, pushing existing rows 568:631 down to 569:632
$ws.Rows("568").Insert()

# Populate the new row 568 with the new weekly data point
$ws.Cells.Item(568, 1).Value = 3
$ws.Cells.Item(568, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(568, 3).Value = "Coquimbo"
$ws.Cells.Item(568, 4).Value = 44946
$ws.Cells.Item(568, 5).Value = 5
$ws.Cells.Item(568, 6).Value = 100112032
$ws.Cells.Item(568, 7).Value = "Zapallo italiano"
$ws.Cells.Item(568, 8).Value = "Sin especificar"
$ws.Cells.Item(568, 9).Value = "Primera"
$ws.Cells.Item(568, 10).Value = 230
$ws.Cells.Item(568, 11).Value = 6000
$ws.Cells.Item(568, 12).Value = 6500
$ws.Cells.Item(568, 13).Value = 6278
$ws.Cells.Item(568, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(568, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(568, 16).Value = 105
$ws.Cells.Item(568, 17).Value = 60
$ws.Cells.Item(568, 18).Value = "Hortaliza"
